# Daily attendance processing - swap the order of "System" and the
# recorded-by email address in column G ("Recorded By") wherever the
# value is exactly "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

$changed = 0
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
        $changed++
    }
}

Write-Host "Updated $changed cell(s) in column G."
